# Updated BMI category name: remove the "BMI_10147" / "Milk" row from
# Sheet1 (row 131). Deleting the entire row shifts every subsequent row
# up by one and Excel automatically renumbers the row references and
# re-points the shared-string indices, which also drops the now-unused
# "BMI_10147" and "Milk" shared strings on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(131).EntireRow.Delete()

# Restore the selection/scroll state recorded in the edited workbook.
$ws.Range("A137").Select()
$excel.ActiveWindow.ScrollRow = 117
$excel.ActiveWindow.ScrollColumn = 1
